# Refresh the crypto symbol list (price + 1h volume change) for the
# "Fri Jan 20 18:54:44 UTC 2023" GitHub Actions data pull.
#
# Price (col D) and Volume(1h) (col E) are stored as plain text in this
# sheet (t="inlineStr"), not numbers/percentages, so each cell's
# NumberFormat is forced to "@" (Text) before the new literal is written -
# this mirrors typing the value into a pre-formatted-as-text cell in
# Excel and keeps e.g. "292.53" / "0.20%" as text instead of being
# reinterpreted as a number or a 0.002-style fraction.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2";  Value = "292.53" },
    @{ Cell = "E2";  Value = "0.20%" },
    @{ Cell = "D3";  Value = "31.19" },
    @{ Cell = "E3";  Value = "0.98%" },
    @{ Cell = "D4";  Value = "4.967" },
    @{ Cell = "E4";  Value = "1.73%" },
    @{ Cell = "D5";  Value = "0.07475" },
    @{ Cell = "E5";  Value = "2.77%" },
    @{ Cell = "D6";  Value = "2.289" },
    @{ Cell = "E6";  Value = "2.74%" },
    @{ Cell = "D7";  Value = "7.775" },
    @{ Cell = "E7";  Value = "1.24%" },
    @{ Cell = "D8";  Value = "0.9192" },
    @{ Cell = "E8";  Value = "2.53%" },
    @{ Cell = "D9";  Value = "0.09455" },
    @{ Cell = "E9";  Value = "19.16%" },
    @{ Cell = "E10"; Value = "4.07%" },
    @{ Cell = "D11"; Value = "0.08327" },
    @{ Cell = "E11"; Value = "2.61%" },
    @{ Cell = "D12"; Value = "0.03270" },
    @{ Cell = "E12"; Value = "5.67%" },
    @{ Cell = "D13"; Value = "0.09907" },
    @{ Cell = "E13"; Value = "-1.34%" },
    @{ Cell = "D14"; Value = "0.001496" },
    @{ Cell = "E14"; Value = "-0.36%" },
    @{ Cell = "D15"; Value = "0.005750" },
    @{ Cell = "E15"; Value = "-1.55%" },
    @{ Cell = "D16"; Value = "3.469" },
    @{ Cell = "E16"; Value = "-0.16%" },
    @{ Cell = "E17"; Value = "1.70%" },
    @{ Cell = "D18"; Value = "2.168" },
    @{ Cell = "E18"; Value = "4.34%" },
    @{ Cell = "E19"; Value = "0.25%" },
    @{ Cell = "E20"; Value = "0.76%" },
    @{ Cell = "E21"; Value = "2.09%" },
    @{ Cell = "D22"; Value = "0.2121" },
    @{ Cell = "E22"; Value = "1.18%" },
    @{ Cell = "D23"; Value = "0.04525" },
    @{ Cell = "E23"; Value = "0.02%" },
    @{ Cell = "D24"; Value = "0.001217" },
    @{ Cell = "E24"; Value = "0.79%" },
    @{ Cell = "D25"; Value = "0.004310" },
    @{ Cell = "E25"; Value = "-7.60%" },
    @{ Cell = "E26"; Value = "0.05%" },
    @{ Cell = "D27"; Value = "0.0003390" },
    @{ Cell = "E27"; Value = "0.05%" },
    @{ Cell = "D39"; Value = "0.01623" },
    @{ Cell = "E39"; Value = "2.98%" },
    @{ Cell = "D40"; Value = "0.04575" },
    @{ Cell = "E40"; Value = "4.03%" },
    @{ Cell = "D41"; Value = "0.007469" },
    @{ Cell = "E41"; Value = "2.20%" },
    @{ Cell = "D42"; Value = "0.009827" },
    @{ Cell = "D43"; Value = "0.1359" },
    @{ Cell = "E43"; Value = "3.32%" },
    @{ Cell = "D44"; Value = "0.002156" },
    @{ Cell = "E44"; Value = "6.82%" },
    @{ Cell = "D45"; Value = "0.009030" },
    @{ Cell = "E45"; Value = "-3.47%" },
    @{ Cell = "D46"; Value = "0.00006094" },
    @{ Cell = "E46"; Value = "6.37%" },
    @{ Cell = "E47"; Value = "0.04%" },
    @{ Cell = "D48"; Value = "2.654" },
    @{ Cell = "E48"; Value = "18.44%" },
    @{ Cell = "D49"; Value = "0.001997" },
    @{ Cell = "E49"; Value = "-30.96%" },
    @{ Cell = "D50"; Value = "0.00002097" },
    @{ Cell = "E50"; Value = "0.04%" },
    @{ Cell = "D51"; Value = "0.0001997" },
    @{ Cell = "E51"; Value = "0.04%" }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $u.Value
}
